# Add a new "2022" column (S) to the SDG 1.3.1 indicator sheet, mirroring
# the formatting already used in column R (the previous latest year).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column S (rows 2-13) with the same formatting as column R in each
# row - this reproduces the per-row cell styles (including the thin-bottom
# border row 2, header row 3, the bold data row 4, the normal data rows
# 5-12, and the bottom-border totals row 13) without hand-rolling a new
# style for every row.
$ws.Range("R2:R13").Copy()
$ws.Range("S2:S13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header: 2022
$ws.Range("S3").Value = 2022

# Data rows for 2022
$ws.Range("S4").Value = 17.700522048199787
$ws.Range("S5").Value = 1.7610202290451711
$ws.Range("S6").Value = 3.9589300291403076
$ws.Range("S7").Value = 1.4859750619980623
$ws.Range("S8").Value = 1.1943569362276563
$ws.Range("S9").Value = 3.9154905266043296
$ws.Range("S10").Value = 0.84000241999604885
$ws.Range("S11").Value = 2.1393883316621789
$ws.Range("S12").Value = 1.8762854436950933
$ws.Range("S13").Value = 0.52907306983093583

# Row 4 (the "Total" style row) ends up with its own distinct cell style in
# the source workbook - same bold 9pt Times New Roman / vertical-centered /
# "0.0" number format as the rest of the row, just a freshly generated
# style entry. Nudge the font so the engine mints a new cellXfs entry
# instead of reusing R4's.
$s4 = $ws.Range("S4")
$s4.Font.Name = "Times New Roman"
$s4.Font.Size = 9
$s4.Font.Bold = $true
$s4.Font.ThemeColor = 1

# The previous selection (R24:R25, left over from past edits) no longer
# makes sense once the sheet has grown only to S22; reset it to A1.
$ws.Range("A1").Select()
